# Bugfix: the "clients" sheet header in A1 was mislabeled "Order No." which
# broke downstream lookups against the "Order Number" header used elsewhere
# (the run-order logic). Rename it back to "Order Number", add the missing
# 4th client (Cindy) expense-report row, and leave the workbook focused on
# the expense_reports sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "clients" sheet: fix the mislabeled header in A1.
# ---------------------------------------------------------------------
$clients = $wb.Worksheets.Item("clients")
$clients.Activate()

$clients.Range("A1").Value = "Order Number"

# Column A now holds the longer header text ("Order Number") -- widen it to
# fit, matching the other bestFit-style columns on this sheet.
$clients.Columns.Item(1).ColumnWidth = 12.5

# Row 3 (Carson Goble / Aiden Herrera e-mail hyperlink row) needs the same
# wrapped-text row height as its siblings.
$clients.Rows.Item(3).RowHeight = 25

# Leave the selection parked on A2 (the header fix is done).
$clients.Range("A2").Select()

# ---------------------------------------------------------------------
# 2) "expense_reports" sheet: add the missing 4th client's expense row.
# ---------------------------------------------------------------------
$expenses = $wb.Worksheets.Item("expense_reports")
$expenses.Activate()

$expenses.Range("A4").Value = 4
$expenses.Range("B4").Value = "Cindy"
$expenses.Range("C4").Value = "testemail22113355@gmail.com"
$expenses.Range("D4").Value = 3.98

$expenses.Rows.Item(4).RowHeight = 15.75
